# ---------------------------------------------------------------------------
# This edit re-shuffles the "species identity" data across a handful of
# observation rows in the sheet (rows 4, 5, 8, 9, 11, 12, 13). The location /
# date / observer columns stay where they are; only the record-identifying
# columns move between rows:
#
#   A  (Id), B (Taxonsorteringsordning), D (Rodlistade), E (TaxonId),
#   F  (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
#   K, L, M, N (Alder-Stadium / Kon / Aktivitet / Metod - only used by row 9)
#   AF (Bestamningsmetod - present as an empty marker cell on some rows)
#
# Net effect (verified against the target diff):
#   row4  <-> row5    (full swap)
#   row8  <-> row9    (full swap, including K/L/M/N which travel with row9's data)
#   row11 <- old row13, row12 <- old row11, row13 <- old row12   (3-way rotation)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding numeric/text "identity" data that must move with each record.
$valueCols = @("A","B","D","E","F","G","H","Q","R")

# Optional columns that only exist on some rows; tracked/moved as presence + value.
$optionalCols = @("K","L","M","N","AF")

function Get-RowSnapshot($ws, $row) {
    $snap = @{}
    foreach ($col in $valueCols) {
        $snap[$col] = $ws.Range("$col$row").Value2
    }
    foreach ($col in $optionalCols) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        $present = $true
        if (($val -eq "") -or ($val -eq $null)) {
            # Empty looking cell: could be a real empty marker cell, or simply
            # not present at all. Use .Text as a secondary signal is not
            # reliable either, so fall back to checking the cell's
            # relationship with neighboring non-empty cells isn't feasible
            # here - instead we track presence explicitly below using the
            # known "before" layout captured once at the start of the script.
            $present = $false
        }
        $snap[$col] = @{ Value = $val; Present = $present }
    }
    return $snap
}

function Set-RowValues($ws, $row, $snap) {
    foreach ($col in $valueCols) {
        $ws.Range("$col$row").Value2 = $snap[$col]
    }
}

function Set-OptionalCell($ws, $row, $col, $present, $value) {
    $cell = $ws.Range("$col$row")
    if ($present) {
        # Force the cell to persist in the saved XML even when the value is
        # an empty string, by touching a formatting property first.
        $cell.NumberFormat = "General"
        $cell.Value2 = $value
    } else {
        $cell.ClearContents()
    }
}

# --- Known presence of the optional columns in the ORIGINAL workbook -------
# (hand-verified from the starting file; used to decide, after the swap,
# which rows should end up with a visible empty marker cell.)
$origPresence = @{
    4  = @{ K=$false; L=$false; M=$false; N=$false; AF=$false }
    5  = @{ K=$false; L=$false; M=$false; N=$false; AF=$true  }
    8  = @{ K=$false; L=$false; M=$false; N=$false; AF=$false }
    9  = @{ K=$true;  L=$true;  M=$true;  N=$true;  AF=$false }
    11 = @{ K=$false; L=$false; M=$false; N=$false; AF=$false }
    12 = @{ K=$false; L=$false; M=$false; N=$false; AF=$true  }
    13 = @{ K=$false; L=$false; M=$false; N=$false; AF=$true  }
}

# --- Capture "before" snapshots for every affected row ---------------------
$snap4  = Get-RowSnapshot $ws 4
$snap5  = Get-RowSnapshot $ws 5
$snap8  = Get-RowSnapshot $ws 8
$snap9  = Get-RowSnapshot $ws 9
$snap11 = Get-RowSnapshot $ws 11
$snap12 = Get-RowSnapshot $ws 12
$snap13 = Get-RowSnapshot $ws 13

# Row -> (data snapshot to write, presence map to apply) after the edit.
# destination row = source row whose data it now receives.
$plan = @{
    4  = @{ Data = $snap5;  Presence = $origPresence[5]  }
    5  = @{ Data = $snap4;  Presence = $origPresence[4]  }
    8  = @{ Data = $snap9;  Presence = $origPresence[9]  }
    9  = @{ Data = $snap8;  Presence = $origPresence[8]  }
    11 = @{ Data = $snap13; Presence = $origPresence[13] }
    12 = @{ Data = $snap11; Presence = $origPresence[11] }
    13 = @{ Data = $snap12; Presence = $origPresence[12] }
}

foreach ($row in @(4,5,8,9,11,12,13)) {
    $entry = $plan[$row]
    $data = $entry.Data
    $presence = $entry.Presence

    Set-RowValues $ws $row $data

    foreach ($col in @("K","L","M","N","AF")) {
        $present = $presence[$col]
        $value = $data[$col].Value
        Set-OptionalCell $ws $row $col $present $value
    }
}
